$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (9th column) to make room for "birthDate"
$ws.Columns.Item(9).Insert()

# Header for the new column
$ws.Cells.Item(1, 9).Value = "birthDate"

# Apply date number format + width to the new column I
$ws.Columns.Item(9).NumberFormat = "m/d/yyyy"
$ws.Columns.Item(9).ColumnWidth = 12.5703125

# Update selection to match target state
$ws.Range("I3").Select()
